# Auto-generated edit script applying the diff's cell-level changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 838.6
$ws.Range("I28").Value = 776.7143
$ws.Range("J28").Value = 983
$ws.Range("K28").Value = 776.7143
$ws.Range("L28").Value = 983
$ws.Range("M28").Value = -291.7143
$ws.Range("N28").Value = -1953
$ws.Range("H51").Value = 6824.75
$ws.Range("I51").Value = 6700
$ws.Range("J51").Value = 6949.5
$ws.Range("K51").Value = 6700
$ws.Range("L51").Value = 6949.5
$ws.Range("M51").Value = -6216
$ws.Range("N51").Value = -7917.5
$ws.Range("H64").Value = 4247.5
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 4247.5
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H74").Value = 3800
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 3800
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860
$ws.Range("H80").Value = 527.5454999999999
$ws.Range("I80").Value = 566.6667
$ws.Range("K80").Value = 1700.0001
$ws.Range("M80").Value = -702.0001
$ws.Range("H83").Value = 527.5454999999999
$ws.Range("I83").Value = 566.6667
$ws.Range("K83").Value = 5100.0003
$ws.Range("M83").Value = -108.0002999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 118
$ws.Range("I5").Value = 118
$ws.Range("K5").Value = 118
$ws.Range("M5").Value = -6
$ws.Range("H61").Value = 28933
$ws.Range("I61").Value = 28933
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 28933
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -28721
$ws.Range("H74").Value = 1435.6364
$ws.Range("I74").Value = 1279.7
$ws.Range("K74").Value = 1279.7
$ws.Range("M74").Value = -405.7
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45676
$ws.Range("H77").Value = 1435.6364
$ws.Range("I77").Value = 1279.7
$ws.Range("K77").Value = 6398.5
$ws.Range("M77").Value = -2030.5
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47340
$ws.Range("H92").Value = 66500
$ws.Range("J92").Value = 66500
$ws.Range("L92").Value = 66500
$ws.Range("N92").Value = -71492
$ws.Range("H136").Value = 28933
$ws.Range("I136").Value = 28933
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 86799
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -84249

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 118
$ws.Range("I4").Value = 118
$ws.Range("K4").Value = 118
$ws.Range("M4").Value = -3
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -127
$ws.Range("H99").Value = 721.75
$ws.Range("J99").Value = 533
$ws.Range("L99").Value = 533
$ws.Range("N99").Value = -3529
$ws.Range("H107").Value = 2873
$ws.Range("I107").Value = 3059.5
$ws.Range("K107").Value = 3059.5
$ws.Range("M107").Value = -1139.5
$ws.Range("H134").Value = 11000
$ws.Range("I134").Value = 10000
$ws.Range("K134").Value = 30000
$ws.Range("M134").Value = -27465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 417.625
$ws.Range("I5").Value = 363
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 363
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = -251
$ws.Range("N5").Value = -1024
$ws.Range("H31").Value = 3439
$ws.Range("I31").Value = 1430
$ws.Range("K31").Value = 1430
$ws.Range("M31").Value = -1135
$ws.Range("H34").Value = 3439
$ws.Range("I34").Value = 1430
$ws.Range("K34").Value = 1430
$ws.Range("M34").Value = -1228
$ws.Range("H41").Value = 1950
$ws.Range("I41").Value = 1950
$ws.Range("K41").Value = 1950
$ws.Range("M41").Value = -1522
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797
$ws.Range("H132").Value = 10679.6
$ws.Range("I132").Value = 8849.75
$ws.Range("K132").Value = 26549.25
$ws.Range("M132").Value = -24019.25
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 501.8
$ws.Range("J121").Value = 661.3333
$ws.Range("L121").Value = 1983.9999
$ws.Range("N121").Value = -4603.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7749.5
$ws.Range("H73").Value = 7749.5
$ws.Range("H80").Value = 18622
$ws.Range("I80").Value = 10902.857
$ws.Range("J80").Value = 36633.332
$ws.Range("K80").Value = 10902.857
$ws.Range("L80").Value = 36633.332
$ws.Range("M80").Value = -9904.857
$ws.Range("N80").Value = -38629.332
$ws.Range("H83").Value = 18622
$ws.Range("I83").Value = 10902.857
$ws.Range("J83").Value = 36633.332
$ws.Range("K83").Value = 54514.285
$ws.Range("L83").Value = 183166.66
$ws.Range("M83").Value = -49522.285
$ws.Range("N83").Value = -193150.66
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = ""
$ws.Range("N103").Value = 0
$ws.Range("H122").Value = 2125
$ws.Range("I122").Value = 2125
$ws.Range("K122").Value = 6375
$ws.Range("M122").Value = -3925
$ws.Range("H126").Value = 3331
$ws.Range("I126").Value = 3331
$ws.Range("K126").Value = 9993
$ws.Range("M126").Value = -7523
$ws.Range("H132").Value = 2092.25
$ws.Range("I132").Value = 1456.3334
$ws.Range("K132").Value = 4369.0002
$ws.Range("M132").Value = -1839.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -888
$ws.Range("H16").Value = 1904.9333
$ws.Range("I16").Value = 2197.2222
$ws.Range("J16").Value = 1466.5
$ws.Range("K16").Value = 2197.2222
$ws.Range("L16").Value = 1466.5
$ws.Range("M16").Value = -2027.2222
$ws.Range("N16").Value = -1806.5
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 5987
$ws.Range("I132").Value = 5987
$ws.Range("K132").Value = 17961
$ws.Range("M132").Value = -15431
$ws.Range("H136").Value = 32800.2
$ws.Range("I136").Value = 32800.2
$ws.Range("K136").Value = 98400.59999999999
$ws.Range("M136").Value = -95850.59999999999
